# Rename the worksheet from "Hoja1" to "Sheet1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Row 1: translate/extend header placeholders
$ws.Range("A1").Value = "Name: {{name}}"
$ws.Range("B1").Value = "Date: {{date}}"

# Row 2: keep existing item placeholders and add new ones
$ws.Range("A2").Value = "ID: [[items.id]]"
$ws.Range("B2").Value = "Qty: [[items.qty]]"
$ws.Range("C2").Value = "Price: [[items.price]]"
$ws.Range("D2").Value = "Item date: [[items.date]]"
$ws.Range("E2").Value = "Missing: [[items.missingProp]]"

# Add a formula cell multiplying quantity by price
$ws.Range("F2").Formula = "=B2*C2"
